$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("constants")

# Row 15 currently holds "start_mdr_introduce_time" -> rename to "mdr_introduce_time"
# (keep its value/comment columns untouched).
$ws.Cells.Item(15, 1).Value = "mdr_introduce_time"

# Row 16 currently holds "end_mdr_introduce_time" (and its related value/comment) -> remove
# the whole row entirely, shifting all subsequent rows up by one.
$ws.Rows.Item(16).Delete()

# Match the saved selection state (active cell moved to A15 after the edit).
$ws.Range("A15").Select() | Out-Null
